# Convert the numeric 0..1 "adcap_score" values in column D (sheet "stressor1")
# into their qualitative labels (none/low/medium/high/NA), add a "lit review"
# note in E7, and refresh the window/selection/zoom state.
#
# New shared strings must be introduced in this exact order so they land at
# the same shared-string indices as the target workbook: none, low, high,
# lit review ("medium" already exists as shared string "medium" and is reused).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stressor1")

# --- category: Spatial Scale of species / EOO (range) ---
$ws.Range("D2").Value = "none"
$ws.Range("D3").Value = "low"
$ws.Range("D4").Value = "medium"
$ws.Range("D5").Value = "high"

# --- category: R (Reproductive Traits) / reproductive strategy ---
$ws.Range("D6").Value = "low"
$ws.Range("D7").Value = "medium"
$ws.Range("E7").Value = "lit review"
$ws.Range("D8").Value = "high"
$ws.Range("D9").Value = "high"

# --- fecundity ---
$ws.Range("D10").Value = "none"
$ws.Range("D11").Value = "none"
$ws.Range("D12").Value = "low"
$ws.Range("D13").Value = "low"
$ws.Range("D14").Value = "low"
$ws.Range("D15").Value = "medium"
$ws.Range("D16").Value = "medium"
$ws.Range("D17").Value = "medium"
$ws.Range("D18").Value = "high"
$ws.Range("D19").Value = "high"

# --- lifetime # reproductive opportunities ---
$ws.Range("D20").Value = "none"
$ws.Range("D21").Value = "low"
$ws.Range("D22").Value = "low"
$ws.Range("D23").Value = "medium"
$ws.Range("D24").Value = "medium"
$ws.Range("D25").Value = "high"

# --- age to 1st reproduction/generation time ---
$ws.Range("D26").Value = "none"
$ws.Range("D27").Value = "low"
$ws.Range("D28").Value = "low"
$ws.Range("D29").Value = "medium"
$ws.Range("D30").Value = "high"

# --- max age ---
$ws.Range("D31").Value = "none"
$ws.Range("D32").Value = "low"
$ws.Range("D33").Value = "low"
$ws.Range("D34").Value = "medium"
$ws.Range("D35").Value = "medium"
$ws.Range("D36").Value = "high"
$ws.Range("D37").Value = "high"

# --- parental investment ---
$ws.Range("D38").Value = "none"
$ws.Range("D39").Value = "high"
$ws.Range("D40").Value = "low"

# --- post-birth/hatching parental dependence ---
$ws.Range("D41").Value = "none"
$ws.Range("D42").Value = "low"
$ws.Range("D43").Value = "medium"
$ws.Range("D44").Value = "medium"
$ws.Range("D45").Value = "high"

# --- global population size ---
$ws.Range("D46").Value = "none"
$ws.Range("D47").Value = "low"
$ws.Range("D48").Value = "medium"
$ws.Range("D49").Value = "high"
$ws.Range("D50").Value = "high"

# --- are there sub-populations? ---
$ws.Range("D51").Value = "high"
$ws.Range("D52").Value = "none"

# --- breeding/nesting range/number of spawning aggregations (fish): number of sites ---
$ws.Range("D53").Value = "none"
$ws.Range("D54").Value = "low"
$ws.Range("D55").Value = "high"
$ws.Range("D56").Value = "high"

# --- sub-population dependence on particular sites ---
$ws.Range("D63").Value = "none"
$ws.Range("D64").Value = "high"

# --- foraging range: number of sites, incl. terrestrial wetlands ---
$ws.Range("D65").Value = "none"
$ws.Range("D66").Value = "low"
$ws.Range("D67").Value = "high"
$ws.Range("D68").Value = "NA"

# --- foraging range: sub-population dependence on particular sites ---
$ws.Range("D75").Value = "none"
$ws.Range("D76").Value = "high"

# --- refresh view state: selection, zoom, and window position ---
$ws.Range("D77").Select()
$excel.ActiveWindow.Zoom = 125

try {
    $win = $wb.Windows.Item(1)
    $win.Left = 2120
    $win.Top = 120
} catch {
}
